# Fruta / hortaliza, semanal
# Insert 4 new weekly rows of data at the top of the "Nectarín" block for
# Feria Lagunitas de Puerto Montt, pushing the existing rows 467-557 down
# to 471-561 and updating the worksheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 467; this shifts the old rows 467:557
# down to 471:561 (Excel copies formatting from the row above, matching
# the existing styling of column D as a date).
$ws.Rows("467:470").Insert()

# Row 467 - Super Queen / Especial
$ws.Range("A467").Value = 4
$ws.Range("B467").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C467").Value = "Los Lagos"
$ws.Range("D467").Value = 44943
$ws.Range("E467").Value = 10
$ws.Range("F467").Value = "Fruta"
$ws.Range("G467").Value = 100103
$ws.Range("H467").Value = "Frutos de hueso (carozo)"
$ws.Range("I467").Value = 100103006
$ws.Range("J467").Value = "Nectarín"
$ws.Range("K467").Value = "Super Queen"
$ws.Range("L467").Value = "Especial"
$ws.Range("M467").Value = 200
$ws.Range("N467").Value = 23000
$ws.Range("O467").Value = 23000
$ws.Range("P467").Value = 23000
$ws.Range("Q467").Value = "$/caja 14 kilos empedrada"
$ws.Range("R467").Value = "Región de O'Higgins"
$ws.Range("S467").Value = 1643
$ws.Range("T467").Value = 14

# Row 468 - Super Queen / Primera
$ws.Range("A468").Value = 4
$ws.Range("B468").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C468").Value = "Los Lagos"
$ws.Range("D468").Value = 44943
$ws.Range("E468").Value = 10
$ws.Range("F468").Value = "Fruta"
$ws.Range("G468").Value = 100103
$ws.Range("H468").Value = "Frutos de hueso (carozo)"
$ws.Range("I468").Value = 100103006
$ws.Range("J468").Value = "Nectarín"
$ws.Range("K468").Value = "Super Queen"
$ws.Range("L468").Value = "Primera"
$ws.Range("M468").Value = 400
$ws.Range("N468").Value = 18000
$ws.Range("O468").Value = 19000
$ws.Range("P468").Value = 18500
$ws.Range("Q468").Value = "$/caja 14 kilos empedrada"
$ws.Range("R468").Value = "Región de O'Higgins"
$ws.Range("S468").Value = 1321
$ws.Range("T468").Value = 14

# Row 469 - Venus / Especial
$ws.Range("A469").Value = 4
$ws.Range("B469").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C469").Value = "Los Lagos"
$ws.Range("D469").Value = 44943
$ws.Range("E469").Value = 10
$ws.Range("F469").Value = "Fruta"
$ws.Range("G469").Value = 100103
$ws.Range("H469").Value = "Frutos de hueso (carozo)"
$ws.Range("I469").Value = 100103006
$ws.Range("J469").Value = "Nectarín"
$ws.Range("K469").Value = "Venus"
$ws.Range("L469").Value = "Especial"
$ws.Range("M469").Value = 200
$ws.Range("N469").Value = 23000
$ws.Range("O469").Value = 23000
$ws.Range("P469").Value = 23000
$ws.Range("Q469").Value = "$/caja 14 kilos empedrada"
$ws.Range("R469").Value = "Región de O'Higgins"
$ws.Range("S469").Value = 1643
$ws.Range("T469").Value = 14

# Row 470 - Venus / Primera
$ws.Range("A470").Value = 4
$ws.Range("B470").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C470").Value = "Los Lagos"
$ws.Range("D470").Value = 44943
$ws.Range("E470").Value = 10
$ws.Range("F470").Value = "Fruta"
$ws.Range("G470").Value = 100103
$ws.Range("H470").Value = "Frutos de hueso (carozo)"
$ws.Range("I470").Value = 100103006
$ws.Range("J470").Value = "Nectarín"
$ws.Range("K470").Value = "Venus"
$ws.Range("L470").Value = "Primera"
$ws.Range("M470").Value = 400
$ws.Range("N470").Value = 18000
$ws.Range("O470").Value = 19000
$ws.Range("P470").Value = 18500
$ws.Range("Q470").Value = "$/caja 14 kilos empedrada"
$ws.Range("R470").Value = "Región de O'Higgins"
$ws.Range("S470").Value = 1321
$ws.Range("T470").Value = 14
